# Update a3.docx via webhook.
#
# The TOC-style heading bookmarks (_Toc...) were regenerated with new
# random suffixes, and the third heading's text was changed from
# "jos jedan samo failed" to "asd".

$d = $word.ActiveDocument

function Rename-Bookmark($oldName, $newName) {
    $bm = $d.Bookmarks.Item($oldName)
    $start = $bm.Start
    $end = $bm.End
    $bm.Delete()
    $rng = $d.Range($start, $end)
    $d.Bookmarks.Add($newName, $rng) | Out-Null
}

# 1) Rename the three auto-generated TOC bookmarks in place.
Rename-Bookmark "_Toc16363637369395480191670182" "_Toc16363774315838753249427539"
Rename-Bookmark "_Toc16363637369749832070168084" "_Toc16363774316166045822362316"
Rename-Bookmark "_Toc16363637370113525278146054" "_Toc16363774316582362598329616"

# 2) Replace the Heading3 text "jos jedan samo failed" with "asd".
#    A plain Range.Text / Find.Execute replace works, but it strips the
#    (empty) run-properties element the original run carried. To keep the
#    run's formatting shape intact, build the replacement text in a throw
#    away paragraph first (so it gets a fresh, empty run-properties of its
#    own), copy that FormattedText onto the heading run, then remove the
#    scratch paragraph again.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$scratchPara = $d.Paragraphs($d.Paragraphs.Count)
$scratchPara.Range.InsertBefore("asd  ")
$templateFt = $scratchPara.Range.FormattedText

$targetBm = $d.Bookmarks.Item("_Toc16363774316582362598329616")
$targetRng = $d.Range($targetBm.Start, $targetBm.End)
$targetRng.FormattedText = $templateFt

$scratchPara2 = $d.Paragraphs($d.Paragraphs.Count)
$scratchRange = $scratchPara2.Range
$delRange = $d.Range($scratchRange.Start - 1, $scratchRange.End)
$delRange.Delete()
